$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46038
$ws.Range("B2").Value = 100.92
$ws.Range("C2").Value = 95.04000000000001
$ws.Range("D2").Value = 94.75
$ws.Range("E2").Value = 94.95
$ws.Range("F2").Value = 94.75
$ws.Range("G2").Value = 95.42
$ws.Range("H2").Value = 103.03
$ws.Range("I2").Value = 121.88
$ws.Range("J2").Value = 133.6
$ws.Range("K2").Value = 130.28
$ws.Range("L2").Value = 118.32
$ws.Range("M2").Value = 111.3
$ws.Range("N2").Value = 109.84
$ws.Range("O2").Value = 108.16
$ws.Range("P2").Value = 116.07
$ws.Range("Q2").Value = 124.38
$ws.Range("R2").Value = 127.64
$ws.Range("S2").Value = 147.46
$ws.Range("T2").Value = 156.89
$ws.Range("U2").Value = 152.79
$ws.Range("V2").Value = 137.01
$ws.Range("W2").Value = 132.02
$ws.Range("X2").Value = 126.75
$ws.Range("Y2").Value = 122.56
$ws.Range("Z2").Value = 118.99
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 146.2
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 154.84
$ws.Range("AE2").Value = "16h-18h"
$ws.Range("AF2").Value = 137.55
$ws.Range("AG2").Value = "0h-14h"
